$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codebook")

# Remove the three rows that no longer apply (Minor_raw, Minor1, Major_minor_inferred),
# which shifts all following rows up by 3 and shrinks the used range from I95 to I92.
$ws.Range("A15:I17").EntireRow.Delete()
